# Advance the payroll week from SEMANA 18 (Apr 29 - May 05 2024) to
# SEMANA 19 (May 06 - May 12 2024), per commit "VIERNES 10 MAYO 2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Main week label (B9) -- downstream cells (H9, B28, H28, B46, H62) are
# formulas that reference this cell (directly or transitively) and will
# recalc automatically.
$ws.Range("B9").Value = "SEMANA  19       DEL    06     Al    12   MAYO     2024"

# Segunda semana de EXTRAS (row 41): 2 x 400 -> 1 x 400
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 400

# Move the active selection/view to match the new snapshot (E42, no pinned
# top-left scroll position).
$null = $ws.Range("E42").Select()
